# Add a new weekly price record for Espinaca (Terminal La Palmera de La Serena)
# at row 226, pushing the existing rows 226:321 down to 227:322.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 226 (shifts rows 226..321 down to 227..322,
# and grows the sheet dimension from A1:R321 to A1:R322).
$ws.Rows("226:226").Insert()

# Populate the newly inserted row with the new observation.
$ws.Range("A226").Value = 8
$ws.Range("B226").Value = "Terminal La Palmera de La Serena"
$ws.Range("C226").Value = "Coquimbo"
$ws.Range("D226").Value = 44837
$ws.Range("E226").Value = 4
$ws.Range("F226").Value = 100112012
$ws.Range("G226").Value = "Espinaca"
$ws.Range("H226").Value = "Sin especificar"
$ws.Range("I226").Value = "Primera"
$ws.Range("J226").Value = 1600
$ws.Range("K226").Value = 450
$ws.Range("L226").Value = 500
$ws.Range("M226").Value = 475
$ws.Range("N226").Value = "`$/atado 300 a 500 gramos"
$ws.Range("O226").Value = "Provincia del Elquí"
$ws.Range("P226").Value = 950
$ws.Range("Q226").Value = 0.5
$ws.Range("R226").Value = "Hortaliza"
